$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.429.10'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '3.418.01'
$ws.Range("E3").Value = '  -3.00%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '584.62'
$ws.Range("E5").Value = '  -2.38%  '
$ws.Range("D6").Value = '137.45'
$ws.Range("E6").Value = '  -4.32%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.417.05'
$ws.Range("E8").Value = '  -2.99%  '
$ws.Range("D9").Value = '0.496'
$ws.Range("E9").Value = '  -0.59%  '
$ws.Range("D10").Value = '7.25'
$ws.Range("E10").Value = '  -6.80%  '
$ws.Range("D11").Value = '0.121'
$ws.Range("E11").Value = '  -10.29%  '
$ws.Range("D12").Value = '0.375'
$ws.Range("E12").Value = '  -7.25%  '
$ws.Range("D13").Value = '3.991.42'
$ws.Range("E13").Value = '  -3.26%  '
$ws.Range("D14").Value = '0.0000179'
$ws.Range("E14").Value = '  -10.28%  '
$ws.Range("E15").Value = '  -1.36%  '
$ws.Range("D16").Value = '3.407.24'
$ws.Range("E16").Value = '  -3.36%  '
$ws.Range("D17").Value = '65.353.33'
$ws.Range("E17").Value = '  -0.72%  '
$ws.Range("D18").Value = '26.11'
$ws.Range("E18").Value = '  -8.80%  '
$ws.Range("E19").Value = '  -10.81%  '
$ws.Range("D20").Value = '5.87'
$ws.Range("E20").Value = '  -5.26%  '
$ws.Range("D21").Value = '13.62'
$ws.Range("E21").Value = '  -5.02%  '
$ws.Range("D22").Value = '384.96'
$ws.Range("E22").Value = '  -6.93%  '
$ws.Range("D23").Value = '0.556'
$ws.Range("E23").Value = '  -7.11%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = '72.54'
$ws.Range("E25").Value = '  -6.00%  '
$ws.Range("D26").Value = '3.551.70'
$ws.Range("E26").Value = '  -3.13%  '
$ws.Range("D27").Value = '0.0000106'
$ws.Range("E27").Value = '  -9.11%  '
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '7.11'
$ws.Range("E29").Value = '  -8.96%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '8.10'
$ws.Range("E30").Value = '  -8.94%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '2.21'
$ws.Range("E31").Value = '  -9.50%  '
$ws.Range("D32").Value = '3.420.63'
$ws.Range("E32").Value = '  -2.85%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '0.143'
$ws.Range("E34").Value = '  -6.28%  '
$ws.Range("D35").Value = '22.78'
$ws.Range("E35").Value = '  -6.41%  '
$ws.Range("D36").Value = '170.47'
$ws.Range("E36").Value = '  -3.09%  '
$ws.Range("D37").Value = '6.78'
$ws.Range("E37").Value = '  -10.01%  '
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").Value = '1.14'
$ws.Range("E38").Value = '  -11.88%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '1.46'
$ws.Range("E39").Value = '  -7.54%  '
$ws.Range("D40").Value = '4.73'
$ws.Range("E40").Value = '  -10.25%  '
$ws.Range("D41").Value = '0.0760'
$ws.Range("E41").Value = '  -7.43%  '
$ws.Range("D42").Value = '0.810'
$ws.Range("E42").Value = '  -5.57%  '
$ws.Range("D43").Value = '43.61'
$ws.Range("E43").Value = '  -3.76%  '
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("D45").Value = '4.39'
$ws.Range("E45").Value = '  -13.33%  '
$ws.Range("D46").Value = '1.60'
$ws.Range("E46").Value = '  -9.84%  '
$ws.Range("D47").Value = '1.10'
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("D48").Value = '22.40'
$ws.Range("E48").Value = '  -1.75%  '
$ws.Range("D49").Value = '6.53'
$ws.Range("E49").Value = '  -7.64%  '
$ws.Range("D50").Value = '2.06'
$ws.Range("E50").Value = '  -14.15%  '
$ws.Range("D51").Value = '2.175.21'
$ws.Range("E51").Value = '  -7.58%  '
